# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-7) ---

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2 (name/count unchanged)
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 82.8

# Row 4: now Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1"
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 77
$ws.Range("D4").Value = 98.5

# Row 5: now Intel(R) Wireless-AC 9560 160MHz - 22.250.0.4
$ws.Range("A5").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.250.0.4"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 98.7

# Row 6: now Intel(R) Wireless-AC 9560 160MHz - 22.250.1.2
$ws.Range("A6").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.250.1.2"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 20

# Row 7: Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2 (name unchanged)
$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 157

# Row 8: Totals
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 271

# --- "Good Drivers" table: Total Samples (column B) updates ---
$ws.Range("B18").Value = 56069
$ws.Range("B19").Value = 449371
$ws.Range("B23").Value = 276086
$ws.Range("B24").Value = 625298
$ws.Range("B29").Value = 453652
$ws.Range("B34").Value = 96091
$ws.Range("B37").Value = 99549
$ws.Range("B38").Value = 77999
$ws.Range("B42").Value = 175767
$ws.Range("B43").Value = 240182
$ws.Range("B49").Value = 684728
$ws.Range("B51").Value = 210188
$ws.Range("B55").Value = 308481
$ws.Range("B60").Value = 443223
$ws.Range("B61").Value = 109665
$ws.Range("B63").Value = 62515
